$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "alleviate" word/definition row (original row 14).
# The runtime will automatically drop the now-unused shared strings
# and shift all subsequent rows up by one.
$ws.Rows("14").Delete()

# Add a picture-filename column (C) that mirrors each word in column A.
$words = @(
  "introspection","philanthropist","antidote","strive","ambidextrous",
  "retrospective","precursors","introvert","gerontocracy","ambiguous",
  "braggart","aggravate","entice","adorn","equilibrium","abhor","connote",
  "endeavor","agile","renovate","curriculum","malevolent","amalgamate","xenophile"
)

for ($i = 0; $i -lt $words.Length; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 3).Value = $words[$i] + ".jpg"
}

# Match the column widths from the updated layout.
$ws.Columns("A").ColumnWidth = 52 - (5/6)
$ws.Columns("B").ColumnWidth = 44 - (5/6)
$ws.Columns("C").ColumnWidth = 32 - (5/6)

# Restore the selection state (active cell A8).
$ws.Range("A8").Select()
